$d = $word.ActiveDocument

# Namespace/package wrapper used for all Range.InsertXML calls below, so the
# inserted fragment stays in the w:body / w:p shape the engine expects while
# replacing only the targeted Range's contents.
function New-RunsXml($innerParagraphXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $innerParagraphXml + '</w:p></w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
}

# ------------------------------------------------------------------
# 1) Drop the stray "_GoBack" bookmark that used to sit right after
#    "This is a change" (first paragraph).
# ------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ------------------------------------------------------------------
# 2) "This is a sub test" -> "This is a sub " + "change" (two runs),
#    only the first occurrence (the sub-bullet right under
#    "This is a change").
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -eq "This is a sub test`r") {
        $full = $d.Range($para.Range.Start, $para.Range.End - 1)
        $full.InsertXML((New-RunsXml('<w:r><w:t xml:space="preserve">This is a sub </w:t></w:r><w:r><w:t>change</w:t></w:r>')))
        break
    }
}

# ------------------------------------------------------------------
# 3) "Added this new item" -> "Added this new " + "change" (two runs),
#    then re-add the "_GoBack" bookmark right after the new text
#    (matches where Word drops it after the most recent edit).
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -eq "Added this new item`r") {
        $full = $d.Range($para.Range.Start, $para.Range.End - 1)
        $full.InsertXML((New-RunsXml('<w:r><w:t xml:space="preserve">Added this new </w:t></w:r><w:r><w:t>change</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>')))
        break
    }
}
